# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Re-sorts / extends the "Estado de Cuenta" detail rows (B16:G26) so the
# data is grouped by period (2110, 2111, 2112, 2201..2206) and the two new
# employees (JESUS ALBERTO MONSALVE MERCADO, HOLVER ANTONIO AMADOR
# PALOMINO) are inserted for their first reported periods (2110 / 2111
# respectively).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: Manuel, period 2206 -> 2110 (first period), value 26650 -> 36341
$ws.Range("E16").Value = "2110"
$ws.Range("F16").Value = 36341

# Row 17: Manuel/2205 -> Jesus Alberto Monsalve Mercado / 2110
$ws.Range("C17").Value = "1044921013"
$ws.Range("D17").Value = "JESUS ALBERTO MONSALVE MERCADO"
$ws.Range("E17").Value = "2110"

# Row 18: Manuel/2204 -> Manuel/2111
$ws.Range("E18").Value = "2111"

# Row 19: Manuel/2203 -> Holver Antonio Amador Palomino / 2111
$ws.Range("C19").Value = "1044922350"
$ws.Range("D19").Value = "HOLVER ANTONIO AMADOR PALOMINO"
$ws.Range("E19").Value = "2111"

# Row 20: Manuel/2202 -> Manuel/2112
$ws.Range("E20").Value = "2112"

# Row 21: Manuel/2201 -> Manuel/2201 (unchanged)

# Row 22: Manuel/2112 -> Manuel/2202
$ws.Range("E22").Value = "2202"

# Row 23: Manuel/2111 -> Manuel/2203
$ws.Range("E23").Value = "2203"

# Row 24: Manuel/2110 -> Manuel/2204
$ws.Range("E24").Value = "2204"

# Row 25: Holver/2111 -> Manuel/2205
$ws.Range("C25").Value = "3883833"
$ws.Range("D25").Value = "MANUEL ANTONIO AMADOR JARAMILLO"
$ws.Range("E25").Value = "2205"

# Row 26: Jesus/2110 -> Manuel/2206, value 36341 -> 26650
$ws.Range("C26").Value = "3883833"
$ws.Range("D26").Value = "MANUEL ANTONIO AMADOR JARAMILLO"
$ws.Range("E26").Value = "2206"
$ws.Range("F26").Value = 26650
